$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates
$ws.Range("B2").Value = 4.9959237662957703
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.0794776331894016
$ws.Range("E2").ClearContents()

# Row 3 updates
$ws.Range("B3").Value = 4.5864924805328933
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 3.0496481571372636
$ws.Range("E3").Value = 7.4335369818962906

# Update selection to match new sqref
$ws.Range("B1:E3").Select()
